# Generate Report for Handoff
#
# The localization pipeline produced a brand-new handoff package for the
# e2e markdown file: a new GUID-based file name and a new xliff content
# hash. The previous handback (target files + handback datetimes) no
# longer applies to the new package, so those fields are cleared back to
# "not yet handed back".

$wb = $excel.ActiveWorkbook

$oldGuid = "b683438f-5d41-41ae-8334-f1fd00e193ac"
$newGuid = "87a785ac-6d15-4bba-9cb8-09a4282cf3d7"
$oldHash = "cce209a9e16e8c4ac9a0f710c55f14bb4dddf846"
$newHash = "37d3c86493fc3b6b53991ba44144af206612f3de"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

$wsOverview.Range("G2").Value = "2016-08-31 09:15:40"

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; HandoffTime = "2016-08-31 09:15:35" },
    @{ Sheet = "de-de"; HandoffTime = "2016-08-31 09:15:40" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Update the source-file hyperlink (A2) to point at the new file name.
    $ws.Range("A2").Value = "$newGuid.md"
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$A$2') {
            $hl.TextToDisplay = "$newGuid.md"
        }
    }

    # Remove the "Latest Target File" hyperlink (I2) entirely - there is
    # no handback target file for the freshly generated package yet.
    $targetLink = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$I$2') {
            $targetLink = $hl
        }
    }
    if ($targetLink -ne $null) {
        $targetLink.Delete()
    }

    # Clear "Latest Target File" and "Latest Handback File" - reset to
    # the not-yet-handed-back state - and drop the hyperlink styling on I2.
    $ws.Range("I2").Value = ""
    $ws.Range("I2").Style = "Normal"
    $ws.Range("J2").Value = ""

    # "Latest Handoff File" now carries the new GUID + content hash.
    $ws.Range("G2").Value = "$newGuid.$newHash.$($lang.Sheet).xlf"

    # "Latest Handoff Datetime" reflects the new handoff.
    $ws.Range("H2").Value = $lang.HandoffTime

    # "Latest Handback DateTime" resets to the null-date sentinel since
    # this package hasn't been handed back yet.
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Column widths for "Latest Target File" / "Latest Handback File"
    # shrink now that they only ever hold short/empty values.
    $ws.Columns.Item(9).ColumnWidth = 17.75
    $ws.Columns.Item(10).ColumnWidth = 20.75
}

Write-Output "Generate Report for Handoff: edits applied"
